# Apply the v0.5.0 "test.xlsx" fixture rewrite:
#  - workbook: date1904 off, sheet renamed hello -> mySheetName
#  - default font bumped to size 12
#  - sheet data replaced with the new sample grid (numbers / booleans / strings / a date)
#  - C3 carries an m/d/yyyy (numFmtId 14) date format

$wb = $excel.ActiveWorkbook
$wb.Date1904 = $false

$ws = $wb.ActiveSheet
$ws.Name = "mySheetName"

# bump the workbook's default font size 11 -> 12 (via the Normal cell style
# so we touch the shared default font rather than stamping a per-range override)
$wb.Styles.Item(1).Font.Size = 12

# clear the old A1:B2 "A1/B1/A2/B2" sample content
$ws.Cells.ClearContents()

# Row 1 - plain numbers
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3

# Row 2 - booleans + a string
$ws.Range("A2").Value = $true
$ws.Range("B2").Value = $false
$ws.Range("D2").Value = "sheetjs"

# Row 3 - strings, a formatted date, and a numeric-looking text value
$ws.Range("A3").Value = "foo"
$ws.Range("B3").Value = "bar"
$ws.Range("C3").Value = 41689.604166666664
$ws.Range("C3").NumberFormat = "m/d/yyyy"
# force "0.3" to stay text instead of being parsed back into the number 0.3
$ws.Range("D3").Value = "'0.3"

# Row 4 - more strings
$ws.Range("A4").Value = "baz"
$ws.Range("C4").Value = "qux"
